$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10k incidents - 2k words")

# 1. Insert 4 new columns before column A, shifting existing A,B,D,E -> E,F,H,I
$ws.Range("A1:D1").EntireColumn.Insert()

# 2. Fill in the new "Results" summary table in columns A:C.
# Text cells are written in the same order the original author typed them
# (this controls the order new entries land in xl/sharedStrings.xml).
$ws.Range("A1").Value = "Results"
$ws.Range("A2").Value = "Incidents analyzed"
$ws.Range("A5").Value = "Common words"
$ws.Range("A4").Value = "Total words analyzed"
$ws.Range("A3").Value = "Common words pool"
$ws.Range("A6").Value = "Cenitex words"
$ws.Range("A7").Value = "Uncommon words"
$ws.Range("A8").Value = "Time"
$ws.Range("B8").Value = "50.31 minutes"

$ws.Range("B2").Value = 10000
$ws.Range("B3").Value = 2000
$ws.Range("B4").Value = 1228224
$ws.Range("B5").Value = 1617
$ws.Range("C5").Formula = "=B5/(SUM(`$B`$5:`$B`$7))"
$ws.Range("B6").Value = 11
$ws.Range("C6").Formula = "=B6/(SUM(`$B`$5:`$B`$7))"
$ws.Range("B7").Value = 33057
$ws.Range("C7").Formula = "=B7/(SUM(`$B`$5:`$B`$7))"

# 3. Merge the title cell across A1:C1
$ws.Range("A1:C1").Merge()

# 4. Formatting
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("A1:C1").HorizontalAlignment = -4108

$ws.Range("B2:B8").NumberFormat = "_-* #,##0_-;-* #,##0_-;_-* ""-""??_-;_-@_-"
$ws.Range("C5:C7").NumberFormat = "0.00%"

# 5. Column widths to match best-fit layout
$ws.Columns.Item(1).ColumnWidth = 19.85546875
$ws.Columns.Item(2).ColumnWidth = 14.85546875
$ws.Columns.Item(3).ColumnWidth = 9.140625

$ws.Range("C8").Select()
